$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "mother_laying_bydate"

# Remove the hidden _xlchart defined names
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# Header row: A1 is (per source) the literal text "MID"; B1:E1 keep their original labels
$ws.Cells.Item(1,1).Value = "MID"
$ws.Cells.Item(1,2).Value = "collect_date"
$ws.Cells.Item(1,3).Value = "n_eggs"
$ws.Cells.Item(1,4).Value = "n_viable"
$ws.Cells.Item(1,5).Value = "pop"

# Data rows 2-60 (A=ID, B=collect_date, C=n_eggs, D=n_viable, E=pop)
# collect_date strings look like dates ("8.17.21") so we force Text format
# first (to stop Excel auto-converting them to date serials), write the
# value, then clear the format back to the workbook default (General/no style)
# so the stored cell has no explicit style, matching the source data.
$ws.Cells.Item(2,1).Value = 119
$bcell = $ws.Cells.Item(2,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.17.21"
$bcell.ClearFormats()
$ws.Cells.Item(2,3).Value = 15
$ws.Cells.Item(2,5).Value = "KL"
$ws.Cells.Item(3,1).Value = 211
$bcell = $ws.Cells.Item(3,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.17.21"
$bcell.ClearFormats()
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,5).Value = "KL"
$ws.Cells.Item(4,1).Value = 92
$bcell = $ws.Cells.Item(4,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.17.21"
$bcell.ClearFormats()
$ws.Cells.Item(4,3).Value = 10
$ws.Cells.Item(4,5).Value = "GV"
$ws.Cells.Item(5,1).Value = 335
$bcell = $ws.Cells.Item(5,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.17.21"
$bcell.ClearFormats()
$ws.Cells.Item(5,3).Value = 16
$ws.Cells.Item(5,5).Value = "PK"
$ws.Cells.Item(6,1).Value = 118
$bcell = $ws.Cells.Item(6,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.17.21"
$bcell.ClearFormats()
$ws.Cells.Item(6,3).Value = 15
$ws.Cells.Item(6,5).Value = "KL"
$ws.Cells.Item(7,1).Value = 114
$bcell = $ws.Cells.Item(7,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.17.21"
$bcell.ClearFormats()
$ws.Cells.Item(7,3).Value = 20
$ws.Cells.Item(7,5).Value = "KL"
$ws.Cells.Item(8,1).Value = 114
$bcell = $ws.Cells.Item(8,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.19.21"
$bcell.ClearFormats()
$ws.Cells.Item(8,3).Value = 22
$ws.Cells.Item(8,5).Value = "KL"
$ws.Cells.Item(9,1).Value = 119
$bcell = $ws.Cells.Item(9,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.19.21"
$bcell.ClearFormats()
$ws.Cells.Item(9,3).Value = 10
$ws.Cells.Item(9,5).Value = "KL"
$ws.Cells.Item(10,1).Value = 16
$bcell = $ws.Cells.Item(10,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.19.21"
$bcell.ClearFormats()
$ws.Cells.Item(10,3).Value = 39
$ws.Cells.Item(10,5).Value = "KL"
$ws.Cells.Item(11,1).Value = 92
$bcell = $ws.Cells.Item(11,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.19.21"
$bcell.ClearFormats()
$ws.Cells.Item(11,3).Value = 20
$ws.Cells.Item(11,5).Value = "GV"
$ws.Cells.Item(12,1).Value = 211
$bcell = $ws.Cells.Item(12,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.19.21"
$bcell.ClearFormats()
$ws.Cells.Item(12,3).Value = 9
$ws.Cells.Item(12,5).Value = "KL"
$ws.Cells.Item(13,1).Value = 118
$bcell = $ws.Cells.Item(13,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.19.21"
$bcell.ClearFormats()
$ws.Cells.Item(13,3).Value = 19
$ws.Cells.Item(13,5).Value = "KL"
$ws.Cells.Item(14,1).Value = 118
$bcell = $ws.Cells.Item(14,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.20.21"
$bcell.ClearFormats()
$ws.Cells.Item(14,3).Value = 1
$ws.Cells.Item(14,5).Value = "KL"
$ws.Cells.Item(15,1).Value = 329
$bcell = $ws.Cells.Item(15,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.20.21"
$bcell.ClearFormats()
$ws.Cells.Item(15,3).Value = 13
$ws.Cells.Item(15,5).Value = "PK"
$ws.Cells.Item(16,1).Value = 363
$bcell = $ws.Cells.Item(16,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.20.21"
$bcell.ClearFormats()
$ws.Cells.Item(16,3).Value = 3
$ws.Cells.Item(16,5).Value = "LP"
$ws.Cells.Item(17,1).Value = 335
$bcell = $ws.Cells.Item(17,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.20.21"
$bcell.ClearFormats()
$ws.Cells.Item(17,3).Value = 3
$ws.Cells.Item(17,5).Value = "PK"
$ws.Cells.Item(18,1).Value = 262
$bcell = $ws.Cells.Item(18,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(18,3).Value = 5
$ws.Cells.Item(18,5).Value = "LP"
$ws.Cells.Item(19,1).Value = 114
$bcell = $ws.Cells.Item(19,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(19,3).Value = 47
$ws.Cells.Item(19,5).Value = "KL"
$ws.Cells.Item(20,1).Value = 164
$bcell = $ws.Cells.Item(20,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(20,3).Value = 3
$ws.Cells.Item(20,5).Value = "LP"
$ws.Cells.Item(21,1).Value = 339
$bcell = $ws.Cells.Item(21,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(21,3).Value = 14
$ws.Cells.Item(21,5).Value = "PK"
$ws.Cells.Item(22,1).Value = 389
$bcell = $ws.Cells.Item(22,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(22,3).Value = 14
$ws.Cells.Item(22,5).Value = "LB"
$ws.Cells.Item(23,1).Value = 367
$bcell = $ws.Cells.Item(23,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(23,3).Value = 8
$ws.Cells.Item(23,5).Value = "LP"
$ws.Cells.Item(24,1).Value = 288
$bcell = $ws.Cells.Item(24,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(24,3).Value = 18
$ws.Cells.Item(24,5).Value = "LB"
$ws.Cells.Item(25,1).Value = 387
$bcell = $ws.Cells.Item(25,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(25,3).Value = 12
$ws.Cells.Item(25,5).Value = "LB"
$ws.Cells.Item(26,1).Value = 349
$bcell = $ws.Cells.Item(26,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(26,3).Value = 24
$ws.Cells.Item(26,5).Value = "HS"
$ws.Cells.Item(27,1).Value = 355
$bcell = $ws.Cells.Item(27,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(27,3).Value = 14
$ws.Cells.Item(27,5).Value = "HS"
$ws.Cells.Item(28,1).Value = 229
$bcell = $ws.Cells.Item(28,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(28,3).Value = 8
$ws.Cells.Item(28,5).Value = "HS"
$ws.Cells.Item(29,1).Value = 28
$bcell = $ws.Cells.Item(29,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(29,3).Value = 25
$ws.Cells.Item(29,5).Value = "HS"
$ws.Cells.Item(30,1).Value = 24
$bcell = $ws.Cells.Item(30,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.25.21"
$bcell.ClearFormats()
$ws.Cells.Item(30,3).Value = 42
$ws.Cells.Item(30,5).Value = "HS"
$ws.Cells.Item(31,1).Value = 119
$bcell = $ws.Cells.Item(31,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(31,3).Value = 10
$ws.Cells.Item(31,5).Value = "KL"
$ws.Cells.Item(32,1).Value = 211
$bcell = $ws.Cells.Item(32,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(32,3).Value = 6
$ws.Cells.Item(32,5).Value = "KL"
$ws.Cells.Item(33,1).Value = 16
$bcell = $ws.Cells.Item(33,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(33,3).Value = 83
$ws.Cells.Item(33,5).Value = "KL"
$ws.Cells.Item(34,1).Value = 118
$bcell = $ws.Cells.Item(34,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(34,3).Value = 40
$ws.Cells.Item(34,5).Value = "KL"
$ws.Cells.Item(35,1).Value = 288
$bcell = $ws.Cells.Item(35,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(35,3).Value = 15
$ws.Cells.Item(35,5).Value = "LB"
$ws.Cells.Item(36,1).Value = 89
$bcell = $ws.Cells.Item(36,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(36,3).Value = 6
$ws.Cells.Item(36,5).Value = "LB"
$ws.Cells.Item(37,1).Value = 298
$bcell = $ws.Cells.Item(37,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(37,3).Value = 10
$ws.Cells.Item(37,5).Value = "LB"
$ws.Cells.Item(38,1).Value = 329
$bcell = $ws.Cells.Item(38,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(38,3).Value = 92
$ws.Cells.Item(38,5).Value = "PK"
$ws.Cells.Item(39,1).Value = 316
$bcell = $ws.Cells.Item(39,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(39,3).Value = 17
$ws.Cells.Item(39,5).Value = "PK"
$ws.Cells.Item(40,1).Value = 2
$bcell = $ws.Cells.Item(40,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(40,3).Value = 19
$ws.Cells.Item(40,5).Value = "PK"
$ws.Cells.Item(41,1).Value = 73
$bcell = $ws.Cells.Item(41,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(41,3).Value = 26
$ws.Cells.Item(41,5).Value = "LW"
$ws.Cells.Item(42,1).Value = 9
$bcell = $ws.Cells.Item(42,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(42,3).Value = 29
$ws.Cells.Item(42,5).Value = "PK"
$ws.Cells.Item(43,1).Value = 317
$bcell = $ws.Cells.Item(43,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(43,3).Value = 16
$ws.Cells.Item(43,5).Value = "PK"
$ws.Cells.Item(44,1).Value = 103
$bcell = $ws.Cells.Item(44,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(44,3).Value = 9
$ws.Cells.Item(44,5).Value = "PK"
$ws.Cells.Item(45,1).Value = 104
$bcell = $ws.Cells.Item(45,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(45,3).Value = 5
$ws.Cells.Item(45,5).Value = "LP"
$ws.Cells.Item(46,1).Value = 363
$bcell = $ws.Cells.Item(46,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(46,3).Value = 11
$ws.Cells.Item(46,5).Value = "LP"
$ws.Cells.Item(47,1).Value = 263
$bcell = $ws.Cells.Item(47,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(47,3).Value = 31
$ws.Cells.Item(47,5).Value = "LP"
$ws.Cells.Item(48,1).Value = 153
$bcell = $ws.Cells.Item(48,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(48,3).Value = 18
$ws.Cells.Item(48,5).Value = "LP"
$ws.Cells.Item(49,1).Value = 247
$bcell = $ws.Cells.Item(49,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(49,3).Value = 12
$ws.Cells.Item(49,5).Value = "LP"
$ws.Cells.Item(50,1).Value = 61
$bcell = $ws.Cells.Item(50,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(50,3).Value = 14
$ws.Cells.Item(50,5).Value = "LP"
$ws.Cells.Item(51,1).Value = 355
$bcell = $ws.Cells.Item(51,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(51,3).Value = 8
$ws.Cells.Item(51,5).Value = "HS"
$ws.Cells.Item(52,1).Value = 24
$bcell = $ws.Cells.Item(52,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(52,3).Value = 10
$ws.Cells.Item(52,5).Value = "HS"
$ws.Cells.Item(53,1).Value = 356
$bcell = $ws.Cells.Item(53,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(53,3).Value = 20
$ws.Cells.Item(53,5).Value = "HS"
$ws.Cells.Item(54,1).Value = 232
$bcell = $ws.Cells.Item(54,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(54,3).Value = 33
$ws.Cells.Item(54,5).Value = "HS"
$ws.Cells.Item(55,1).Value = 33
$bcell = $ws.Cells.Item(55,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(55,3).Value = 35
$ws.Cells.Item(55,5).Value = "HS"
$ws.Cells.Item(56,1).Value = 348
$bcell = $ws.Cells.Item(56,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(56,3).Value = 22
$ws.Cells.Item(56,5).Value = "HS"
$ws.Cells.Item(57,1).Value = 129
$bcell = $ws.Cells.Item(57,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(57,3).Value = 21
$ws.Cells.Item(57,5).Value = "HS"
$ws.Cells.Item(58,1).Value = 41
$bcell = $ws.Cells.Item(58,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(58,3).Value = 44
$ws.Cells.Item(58,5).Value = "HS"
$ws.Cells.Item(59,1).Value = 29
$bcell = $ws.Cells.Item(59,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(59,3).Value = 14
$ws.Cells.Item(59,5).Value = "HS"
$ws.Cells.Item(60,1).Value = 22
$bcell = $ws.Cells.Item(60,2)
$bcell.NumberFormat = "@"
$bcell.Value = "8.27.21"
$bcell.ClearFormats()
$ws.Cells.Item(60,3).Value = 5
$ws.Cells.Item(60,5).Value = "HS"
